$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.612.78'
Set-TextValue 'E2' '  +0.55%  '
Set-TextValue 'D3' '1.744.42'
Set-TextValue 'E3' '  +0.91%  '
Set-TextValue 'D4' '0.9994'
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '247.05'
Set-TextValue 'E5' '  +1.25%  '
Set-TextValue 'D6' '0.9996'
Set-TextValue 'D7' '0.4931'
Set-TextValue 'E7' '  +2.47%  '
Set-TextValue 'D8' '0.2683'
Set-TextValue 'E8' '  +0.17%  '
Set-TextValue 'D9' '0.06291'
Set-TextValue 'E9' '  +0.92%  '
Set-TextValue 'D10' '1.741.70'
Set-TextValue 'E10' '  +0.68%  '
Set-TextValue 'D11' '0.07058'
Set-TextValue 'E11' '  -1.11%  '
Set-TextValue 'D12' '15.79'
Set-TextValue 'E12' '  +0.27%  '
Set-TextValue 'D13' '0.6174'
Set-TextValue 'E13' '  -0.23%  '
Set-TextValue 'D14' '4.599'
Set-TextValue 'E14' '  +0.94%  '
Set-TextValue 'E15' '  +1.50%  '
Set-TextValue 'D16' '0.9997'
Set-TextValue 'E16' '  -0.02%  '
Set-TextValue 'D17' '26.626.72'
Set-TextValue 'E17' '  +0.55%  '
Set-TextValue 'D18' '0.000007314'
Set-TextValue 'E18' '  +5.09%  '
Set-TextValue 'D19' '0.9997'
Set-TextValue 'E19' '  +0.00%  '
Set-TextValue 'D20' '11.59'
Set-TextValue 'E20' '  -1.38%  '
Set-TextValue 'D21' '1.961.85'
Set-TextValue 'E21' '  +0.39%  '
Set-TextValue 'D22' '4.599'
Set-TextValue 'E22' '  +0.86%  '
Set-TextValue 'D23' '8.747'
Set-TextValue 'E23' '  -2.15%  '
Set-TextValue 'D24' '5.281'
Set-TextValue 'E24' '  -0.85%  '
Set-TextValue 'D25' '139.69'
Set-TextValue 'E25' '  +2.29%  '
Set-TextValue 'D26' '15.50'
Set-TextValue 'D27' '1.424'
Set-TextValue 'E27' '  +1.23%  '
Set-TextValue 'B28' 'LidoDAOToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D28' '1.771'
Set-TextValue 'E28' '  -1.46%  '
Set-TextValue 'B29' 'BitcoinCash'
Set-TextValue 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D29' '107.84'
Set-TextValue 'E29' '  +1.02%  '
Set-TextValue 'D30' '4.060'
Set-TextValue 'E30' '  +1.76%  '
Set-TextValue 'D31' '0.08060'
Set-TextValue 'E31' '  +0.40%  '
Set-TextValue 'D32' '3.748'
Set-TextValue 'E32' '  +0.13%  '
Set-TextValue 'D33' '0.04626'
Set-TextValue 'E33' '  +1.42%  '
Set-TextValue 'D35' '2.609'
Set-TextValue 'E35' '  -0.22%  '
Set-TextValue 'E36' '  +2.89%  '
Set-TextValue 'D37' '0.6402'
Set-TextValue 'E37' '  -0.23%  '
Set-TextValue 'D38' '2.090'
Set-TextValue 'E38' '  +4.69%  '
Set-TextValue 'D39' '0.9010'
Set-TextValue 'E39' '  -4.71%  '
Set-TextValue 'D40' '2.423'
Set-TextValue 'E40' '  +0.42%  '
Set-TextValue 'E41' '  -0.30%  '
Set-TextValue 'D42' '0.01507'
Set-TextValue 'E42' '  +0.18%  '
Set-TextValue 'D43' '101.82'
Set-TextValue 'E43' '  -5.61%  '
Set-TextValue 'D44' '5.431'
Set-TextValue 'E44' '  -4.35%  '
Set-TextValue 'D45' '0.3937'
Set-TextValue 'E45' '  +0.38%  '
Set-TextValue 'D46' '6.917'
Set-TextValue 'E46' '  -1.14%  '
Set-TextValue 'D47' '0.1188'
Set-TextValue 'E47' '  -0.62%  '
Set-TextValue 'D48' '0.05399'
Set-TextValue 'E48' '  +1.41%  '
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '7.871'
Set-TextValue 'E49' '  -0.52%  '
Set-TextValue 'B50' 'Elrond'
Set-TextValue 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D50' '30.62'
Set-TextValue 'E50' '  -0.95%  '
Set-TextValue 'D51' '1.268'
Set-TextValue 'E51' '  -0.63%  '

Write-Host "Applied all cell updates"
